# Insert a new weekly data row at row 21 (pushing existing rows 21-28 down to 22-29)
# and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 21; this shifts rows 21-28 -> 22-29
$ws.Rows.Item(21).Insert()

# Copy the date cell style (numeric/date format) from the row below (now row 22, col D)
# into the newly inserted row 21's D cell, so the date renders correctly.
$ws.Cells.Item(22, 4).Copy()
$ws.Cells.Item(21, 4).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row 21 values
$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(21, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(21, 4).Value = 45240
$ws.Cells.Item(21, 5).Value = 15
$ws.Cells.Item(21, 6).Value = "Fruta"
$ws.Cells.Item(21, 7).Value = 100107
$ws.Cells.Item(21, 8).Value = "Otros"
$ws.Cells.Item(21, 9).Value = 100107002
$ws.Cells.Item(21, 10).Value = "Chirimoya"
$ws.Cells.Item(21, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(21, 12).Value = "Primera"
$ws.Cells.Item(21, 13).Value = 225
$ws.Cells.Item(21, 14).Value = 18000
$ws.Cells.Item(21, 15).Value = 20000
$ws.Cells.Item(21, 16).Value = 18889
$ws.Cells.Item(21, 17).Value = "$/caja 13 kilos"
$ws.Cells.Item(21, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(21, 19).Value = 1453
$ws.Cells.Item(21, 20).Value = 13

$wb.Save()
